$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (1-based column index -> new width).
# Note: Excel's ColumnWidth property (character units) differs from the
# raw OOXML <col width="..."> value stored on disk by a constant offset of
# 5/6 (0.8333...) for this workbook's font metrics. To land on an exact
# target stored width, subtract that offset before assigning ColumnWidth.
$offset = 5/6
$ws.Columns.Item(1).ColumnWidth = 81 - $offset
$ws.Columns.Item(2).ColumnWidth = 35 - $offset
$ws.Columns.Item(4).ColumnWidth = 40 - $offset
$ws.Columns.Item(5).ColumnWidth = 42 - $offset
$ws.Columns.Item(6).ColumnWidth = 30 - $offset
$ws.Columns.Item(7).ColumnWidth = 32 - $offset
$ws.Columns.Item(8).ColumnWidth = 27 - $offset
$ws.Columns.Item(9).ColumnWidth = 29 - $offset
$ws.Columns.Item(10).ColumnWidth = 35 - $offset
$ws.Columns.Item(11).ColumnWidth = 37 - $offset
$ws.Columns.Item(12).ColumnWidth = 30 - $offset
$ws.Columns.Item(13).ColumnWidth = 32 - $offset

# Update header row text values (row 1)
$ws.Range("B1").Value = "div_testRuns_internalRoleCellName"
$ws.Range("D1").Value = "link_projectLinks_internalRoleLinkName"
$ws.Range("E1").Value = "link_projectLinks_internalRoleLinkName_1"
$ws.Range("F1").Value = "link_projectLinks_project_id"
$ws.Range("G1").Value = "link_projectLinks_project_id_1"
$ws.Range("H1").Value = "link_projectLinks_team_id"
$ws.Range("I1").Value = "link_projectLinks_team_id_1"
$ws.Range("J1").Value = "link_projectLinks_test_project_id"
$ws.Range("K1").Value = "link_projectLinks_test_project_id_1"
$ws.Range("L1").Value = "link_projectLinks_trNthChild"
$ws.Range("M1").Value = "link_projectLinks_trNthChild_1"

# Update data row (row 2)
$ws.Range("A2").Value = "Data Files/AI-Generated/Common/scheduleTestRunAndConfigureEnvironment-test-data"
